# fixed harvester column in rnasamples -- holly added S.GISH to harvester
# in bioSamples
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B22").Value = "S.GISH"

$ws.Range("B:B").Select()
